$wb = $excel.ActiveWorkbook

# Sheet "z1,1" (index 1): turn off the two 1-markers in rows 5 and 6
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("H5").Value = 0
$ws1.Range("E6").Value = 0

# Sheet "z2,1" (index 2): row 6 marker moves from B6 to C6, and E6 also becomes 1
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B6").Value = 0
$ws2.Range("C6").Value = 1
$ws2.Range("E6").Value = 1

# Sheet "z1,3" (index 5): clear F1, C6, F10 and set I4
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("F1").Value = 0
$ws5.Range("I4").Value = 1
$ws5.Range("C6").Value = 0
$ws5.Range("F10").Value = 0

# Sheet "z2,3" (index 6): set F1, H5, B6, F10 and clear I4
$ws6 = $wb.Worksheets.Item(6)
$ws6.Range("F1").Value = 1
$ws6.Range("I4").Value = 0
$ws6.Range("H5").Value = 1
$ws6.Range("B6").Value = 1
$ws6.Range("F10").Value = 1
